# Apply updated Fitness values (column C) based on the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 13379
    3 = 11228
    4 = 10447
    5 = 10447
    6 = 10447
    7 = 10447
    8 = 10447
    9 = 10306
    10 = 9532
    11 = 9532
    12 = 9532
    13 = 9532
    14 = 9532
    15 = 9420
    16 = 8656
    17 = 8656
    18 = 8595
    19 = 8595
    20 = 8595
    21 = 8595
    22 = 8595
    23 = 8595
    24 = 8595
    25 = 8595
    26 = 8595
    27 = 8303
    28 = 8303
    29 = 8303
    30 = 8303
    31 = 8303
    32 = 8303
    33 = 8303
    34 = 8303
    35 = 8303
    36 = 8303
    37 = 8283
    38 = 8283
    39 = 8233
    40 = 8233
    41 = 8233
    42 = 7898
    43 = 7898
    44 = 7872
    45 = 7872
    46 = 7872
    47 = 7870
    48 = 7870
    49 = 7870
    50 = 7870
    51 = 7870
    52 = 7870
    53 = 7870
    54 = 7870
    55 = 7870
    56 = 7870
    57 = 7870
    58 = 7870
    59 = 7870
    60 = 7870
    61 = 7870
    62 = 7870
    63 = 7870
    64 = 7870
    65 = 7870
    66 = 7870
    67 = 7870
    68 = 7870
    69 = 7870
    70 = 7870
    71 = 7870
    72 = 7870
    73 = 7870
    74 = 7870
    75 = 7870
    76 = 7870
    77 = 7870
    78 = 7870
    79 = 7870
    80 = 7870
    81 = 7870
    82 = 7870
    83 = 7870
    84 = 7870
    85 = 7293
    86 = 7293
    87 = 7293
    88 = 7293
    89 = 7293
    90 = 7293
    91 = 7293
    92 = 7293
    93 = 7293
    94 = 7293
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
